$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Search turned up rows whose Agent/Name (column B) only matched a stale
# "New EP User" placeholder and rows that had not actually been processed
# yet (columns B "Name" and E "Estado de RITM"). Clear those out.
$ws.Range("B5:B54").ClearContents()
$ws.Range("E5:E54").ClearContents()

# The bottom rows that were tagged with the placeholder "New EP User"
# actually belong to Juan Carlos DIOSES - fix the agent name.
$ws.Range("B76:B80").Value = "Juan Carlos DIOSES"

# Leave the selection where the review ended.
$ws.Range("C9").Select() | Out-Null
